$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value without Excel reinterpreting
# numeric-looking strings (e.g. "1.00", "68.118.13") as numbers, and without
# leaving a permanent text-format style applied to the cell afterwards.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.118.13"
Set-TextValue $ws.Range("E2") "  +1.60%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.587.42"
Set-TextValue $ws.Range("E3") "  +0.04%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.12%  "

# Row 5
Set-TextValue $ws.Range("D5") "209.18"
Set-TextValue $ws.Range("E5") "  +8.90%  "

# Row 6
Set-TextValue $ws.Range("D6") "571.13"
Set-TextValue $ws.Range("E6") "  -0.43%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.612"
Set-TextValue $ws.Range("E7") "  -0.85%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.21%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.680"
Set-TextValue $ws.Range("E9") "  +0.34%  "

# Row 10
Set-TextValue $ws.Range("D10") "62.43"
Set-TextValue $ws.Range("E10") "  +12.38%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.148"
Set-TextValue $ws.Range("E11") "  -1.74%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000281"
Set-TextValue $ws.Range("E12") "  +4.56%  "

# Row 13
Set-TextValue $ws.Range("D13") "10.30"
Set-TextValue $ws.Range("E13") "  +4.63%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.150.19"
Set-TextValue $ws.Range("E14") "  -0.05%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.576.31"
Set-TextValue $ws.Range("E15") "  -0.17%  "

# Row 16
Set-TextValue $ws.Range("E16") "  +0.65%  "

# Row 17
Set-TextValue $ws.Range("D17") "19.20"
Set-TextValue $ws.Range("E17") "  +4.46%  "

# Row 18
Set-TextValue $ws.Range("D18") "67.894.60"
Set-TextValue $ws.Range("E18") "  +1.45%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.21"
Set-TextValue $ws.Range("E19") "  +0.54%  "

# Row 20
Set-TextValue $ws.Range("E20") "  +0.08%  "

# Row 21
Set-TextValue $ws.Range("D21") "402.97"
Set-TextValue $ws.Range("E21") "  +1.20%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.16"
Set-TextValue $ws.Range("E22") "  -0.63%  "

# Row 23
Set-TextValue $ws.Range("D23") "12.45"
Set-TextValue $ws.Range("E23") "  +10.07%  "

# Row 24
Set-TextValue $ws.Range("D24") "84.56"
Set-TextValue $ws.Range("E24") "  -1.36%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.89"
Set-TextValue $ws.Range("E25") "  -2.03%  "

# Row 26
Set-TextValue $ws.Range("D26") "3.91"
Set-TextValue $ws.Range("E26") "  +7.76%  "

# Row 27
Set-TextValue $ws.Range("D27") "12.48"
Set-TextValue $ws.Range("E27") "  +0.17%  "

# Row 28
Set-TextValue $ws.Range("D28") "9.26"
Set-TextValue $ws.Range("E28") "  +3.60%  "

# Row 29
Set-TextValue $ws.Range("D29") "7.65"
Set-TextValue $ws.Range("E29") "  -0.09%  "

# Row 30
Set-TextValue $ws.Range("D30") "31.55"
Set-TextValue $ws.Range("E30") "  +1.24%  "

# Row 31
Set-TextValue $ws.Range("D31") "681.52"
Set-TextValue $ws.Range("E31") "  +7.86%  "

# Row 32
Set-TextValue $ws.Range("D32") "12.11"
Set-TextValue $ws.Range("E32") "  -0.56%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.113"
Set-TextValue $ws.Range("E33") "  -1.51%  "

# Row 34
Set-TextValue $ws.Range("D34") "63.22"
Set-TextValue $ws.Range("E34") "  -1.01%  "

# Row 35
Set-TextValue $ws.Range("D35") "40.98"
Set-TextValue $ws.Range("E35") "  -2.51%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.410"
Set-TextValue $ws.Range("E36") "  +2.77%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.19%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.21"
Set-TextValue $ws.Range("E38") "  +8.30%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0₃0750"
Set-TextValue $ws.Range("E39") "  -1.57%  "

# Row 40
Set-TextValue $ws.Range("D40") "3.164.56"
Set-TextValue $ws.Range("E40") "  +0.19%  "

# Row 41
Set-TextValue $ws.Range("E41") "  +0.17%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D42") "3.05"
Set-TextValue $ws.Range("E42") "  +20.00%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D43") "0.997"
Set-TextValue $ws.Range("E43") "  -0.08%  "

# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D44") "2.66"
Set-TextValue $ws.Range("E44") "  -1.87%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D45") "0.0412"
Set-TextValue $ws.Range("E45") "  -0.64%  "

# Row 46
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D46") "2.74"
Set-TextValue $ws.Range("E46") "  +8.47%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.12"
Set-TextValue $ws.Range("E47") "  +1.64%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -0.13%  "

# Row 49
Set-TextValue $ws.Range("D49") "8.68"
Set-TextValue $ws.Range("E49") "  +2.04%  "

# Row 50
Set-TextValue $ws.Range("D50") "138.21"
Set-TextValue $ws.Range("E50") "  -1.50%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.69"
Set-TextValue $ws.Range("E51") "  -1.91%  "

